$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 91810
$ws.Range("B4").Value = 92464
$ws.Range("B5").Value = 91867
$ws.Range("B6").Value = 92298
